# Apply updated TPM values to LR-pairs sheet (Tgfa-Erbb2)
# Rows 2-7 get refreshed numeric values (and a few re-labelled target/sending
# clusters), and four new rows (8-11) are appended for the two new clusters
# "Inflammatory-Mac" and "Resolving-Mac".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: refreshed TPM-derived values
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5909176666666667
$ws.Range("H2").Value = 1.772753
$ws.Range("I2").Value = 0.9937758428931484
$ws.Range("J2").Value = 0.9937758428931482
$ws.Range("M2").Value = 3.248861
$ws.Range("N2").Value = 9.746583
$ws.Range("O2").Value = 0.3767295193213461
$ws.Range("P2").Value = 0.376729519321346
$ws.Range("Q2").Value = 1.919809361444333
$ws.Range("R2").Value = 17.278284252999
$ws.Range("S2").Value = 0.3743846956063014
$ws.Range("T2").Value = 0.3743846956063012

# Row 3: refreshed TPM-derived values
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5909176666666667
$ws.Range("H3").Value = 1.772753
$ws.Range("I3").Value = 0.9937758428931484
$ws.Range("J3").Value = 0.9937758428931482
$ws.Range("M3").Value = 3.599001333333334
$ws.Range("O3").Value = 0.417330886838049
$ws.Range("P3").Value = 0.4173308868380489
$ws.Range("Q3").Value = 2.126713470223556
$ws.Range("R3").Value = 19.140421232012
$ws.Range("S3").Value = 0.4147333538328272
$ws.Range("T3").Value = 0.4147333538328271

# Row 4: refreshed TPM-derived values
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5909176666666667
$ws.Range("H4").Value = 1.772753
$ws.Range("I4").Value = 0.9937758428931484
$ws.Range("J4").Value = 0.9937758428931482
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02714433333333334
$ws.Range("N4").Value = 0.081433
$ws.Range("O4").Value = 0.003147586692371591
$ws.Range("P4").Value = 0.003147586692371591
$ws.Range("Q4").Value = 0.01604006611655556
$ws.Range("R4").Value = 0.144360595049
$ws.Range("S4").Value = 0.003127995618290835
$ws.Range("T4").Value = 0.003127995618290834

# Row 5: refreshed TPM-derived values
$ws.Range("A5").Value = "ECs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.5909176666666667
$ws.Range("H5").Value = 1.772753
$ws.Range("I5").Value = 0.9937758428931484
$ws.Range("J5").Value = 0.9937758428931482
$ws.Range("M5").Value = 1.711959
$ws.Range("N5").Value = 5.135877
$ws.Range("O5").Value = 0.1985143381535413
$ws.Range("P5").Value = 0.1985143381535413
$ws.Range("Q5").Value = 1.011626817709
$ws.Range("R5").Value = 9.104641359381
$ws.Range("S5").Value = 0.197278753724911
$ws.Range("T5").Value = 0.1972787537249109

# Row 6: refreshed TPM-derived values
$ws.Range("A6").Value = "ECs"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.5909176666666667
$ws.Range("H6").Value = 1.772753
$ws.Range("I6").Value = 0.9937758428931484
$ws.Range("J6").Value = 0.9937758428931482
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03689
$ws.Range("N6").Value = 0.11067
$ws.Range("O6").Value = 0.004277668994692127
$ws.Range("P6").Value = 0.004277668994692126
$ws.Range("Q6").Value = 0.02179895272333333
$ws.Range("R6").Value = 0.19619057451
$ws.Range("S6").Value = 0.004251044110818055
$ws.Range("T6").Value = 0.004251044110818053

# Row 7: refreshed TPM-derived values
$ws.Range("D7").Value = "ECs"
$ws.Range("I7").Value = 0.006224157106851674
$ws.Range("J7").Value = 0.006224157106851673
$ws.Range("M7").Value = 3.248861
$ws.Range("N7").Value = 9.746583
$ws.Range("O7").Value = 0.3767295193213461
$ws.Range("P7").Value = 0.376729519321346
$ws.Range("Q7").Value = 0.012024034561
$ws.Range("R7").Value = 0.108216311049
$ws.Range("S7").Value = 0.002344823715044772
$ws.Range("T7").Value = 0.002344823715044771

# Row 8: new row
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Tgfa"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.003701
$ws.Range("H8").Value = 0.011103
$ws.Range("I8").Value = 0.006224157106851674
$ws.Range("J8").Value = 0.006224157106851673
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.599001333333334
$ws.Range("N8").Value = 10.797004
$ws.Range("O8").Value = 0.417330886838049
$ws.Range("P8").Value = 0.4173308868380489
$ws.Range("Q8").Value = 0.01331990393466667
$ws.Range("R8").Value = 0.119879135412
$ws.Range("S8").Value = 0.002597533005221754
$ws.Range("T8").Value = 0.002597533005221754

# Row 9: new row
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Tgfa"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.003701
$ws.Range("H9").Value = 0.011103
$ws.Range("I9").Value = 0.006224157106851674
$ws.Range("J9").Value = 0.006224157106851673
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.02714433333333334
$ws.Range("N9").Value = 0.081433
$ws.Range("O9").Value = 0.003147586692371591
$ws.Range("P9").Value = 0.003147586692371591
$ws.Range("Q9").Value = 0.0001004611776666667
$ws.Range("R9").Value = 0.000904150599
$ws.Range("S9").Value = 0.00001959107408075639
$ws.Range("T9").Value = 0.00001959107408075639

# Row 10: new row
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Tgfa"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.003701
$ws.Range("H10").Value = 0.011103
$ws.Range("I10").Value = 0.006224157106851674
$ws.Range("J10").Value = 0.006224157106851673
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.711959
$ws.Range("N10").Value = 5.135877
$ws.Range("O10").Value = 0.1985143381535413
$ws.Range("P10").Value = 0.1985143381535413
$ws.Range("Q10").Value = 0.006335960259
$ws.Range("R10").Value = 0.057023642331
$ws.Range("S10").Value = 0.001235584428630321
$ws.Range("T10").Value = 0.00123558442863032

# Row 11: new row
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Tgfa"
$ws.Range("C11").Value = "Erbb2"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.003701
$ws.Range("H11").Value = 0.011103
$ws.Range("I11").Value = 0.006224157106851674
$ws.Range("J11").Value = 0.006224157106851673
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.03689
$ws.Range("N11").Value = 0.11067
$ws.Range("O11").Value = 0.004277668994692127
$ws.Range("P11").Value = 0.004277668994692126
$ws.Range("Q11").Value = 0.00013652989
$ws.Range("R11").Value = 0.00122876901
$ws.Range("S11").Value = 0.00002662488387407206
$ws.Range("T11").Value = 0.00002662488387407205

